# Auto-generated edit script applying numeric corrections to the
# Leve profit calculation sheets (currentAveragePrice / Leve price /
# profit columns H-N) across ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

$ws.Range("H137").Value = 7877.4443
$ws.Range("J137").Value = 15000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -50100

$ws.Range("H138").Value = 3859.2727
$ws.Range("J138").Value = 4931.75
$ws.Range("L138").Value = 14795.25
$ws.Range("N138").Value = -25075.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1556.0952
$ws.Range("I32").Value = 1556.0952
$ws.Range("K32").Value = 1556.0952
$ws.Range("M32").Value = -1269.0952

$ws.Range("H61").Value = 3240.7273
$ws.Range("J61").Value = 4548.3076
$ws.Range("L61").Value = 4548.3076
$ws.Range("N61").Value = -4972.3076

$ws.Range("H97").Value = 2333.875
$ws.Range("I97").Value = 2036.4
$ws.Range("K97").Value = 2036.4
$ws.Range("M97").Value = -1540.4

$ws.Range("H132").Value = 2193.9666
$ws.Range("I132").Value = 994.4211
$ws.Range("J132").Value = 4265.909
$ws.Range("K132").Value = 2983.2633
$ws.Range("L132").Value = 12797.727
$ws.Range("M132").Value = -453.2633000000001
$ws.Range("N132").Value = -17857.727

$ws.Range("H136").Value = 3240.7273
$ws.Range("J136").Value = 4548.3076
$ws.Range("L136").Value = 13644.9228
$ws.Range("N136").Value = -18744.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1328.6666
$ws.Range("I80").Value = 94.8
$ws.Range("J80").Value = 2210
$ws.Range("K80").Value = 94.8
$ws.Range("L80").Value = 2210
$ws.Range("M80").Value = 903.2
$ws.Range("N80").Value = -4206

$ws.Range("H83").Value = 1328.6666
$ws.Range("I83").Value = 94.8
$ws.Range("J83").Value = 2210
$ws.Range("K83").Value = 474
$ws.Range("L83").Value = 11050
$ws.Range("M83").Value = 4518
$ws.Range("N83").Value = -21034

$ws.Range("H86").Value = 1353.5
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 1707
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 1707
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -3953

$ws.Range("H89").Value = 1353.5
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 1707
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 8535
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -19767

$ws.Range("H94").Value = 932.6
$ws.Range("I94").Value = 715.1429000000001
$ws.Range("K94").Value = 715.1429000000001
$ws.Range("M94").Value = -264.1429000000001

$ws.Range("H99").Value = 1955
$ws.Range("I99").Value = 1955
$ws.Range("K99").Value = 1955
$ws.Range("M99").Value = -457

$ws.Range("H107").Value = 658.25
$ws.Range("J107").Value = 756.5
$ws.Range("L107").Value = 756.5
$ws.Range("N107").Value = -4596.5

$ws.Range("H134").Value = 2542.5881
$ws.Range("I134").Value = 1326.5
$ws.Range("J134").Value = 3623.5557
$ws.Range("K134").Value = 3979.5
$ws.Range("L134").Value = 10870.6671
$ws.Range("M134").Value = -1444.5
$ws.Range("N134").Value = -15940.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2863.8948
$ws.Range("I31").Value = 1563.8889
$ws.Range("J31").Value = 4033.9
$ws.Range("K31").Value = 1563.8889
$ws.Range("L31").Value = 4033.9
$ws.Range("M31").Value = -1268.8889
$ws.Range("N31").Value = -4623.9

$ws.Range("H34").Value = 2863.8948
$ws.Range("I34").Value = 1563.8889
$ws.Range("J34").Value = 4033.9
$ws.Range("K34").Value = 1563.8889
$ws.Range("L34").Value = 4033.9
$ws.Range("M34").Value = -1361.8889
$ws.Range("N34").Value = -4437.9

$ws.Range("H68").Value = 54647.5
$ws.Range("J68").Value = 79295
$ws.Range("L68").Value = 79295
$ws.Range("N68").Value = -80793

$ws.Range("H71").Value = 54647.5
$ws.Range("J71").Value = 79295
$ws.Range("L71").Value = 237885
$ws.Range("N71").Value = -245373

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 142.66667
$ws.Range("J12").Value = 169.5
$ws.Range("L12").Value = 508.5
$ws.Range("N12").Value = -854.5

$ws.Range("H38").Value = 43.11111
$ws.Range("I38").Value = 44
$ws.Range("J38").Value = 36
$ws.Range("K38").Value = 132
$ws.Range("L38").Value = 108
$ws.Range("M38").Value = 215
$ws.Range("N38").Value = -802

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

$ws.Range("H107").Value = 1000
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840

$ws.Range("H113").Value = 958
$ws.Range("I113").Value = 403
$ws.Range("K113").Value = 1209
$ws.Range("M113").Value = 961

$ws.Range("H122").Value = 848.8
$ws.Range("J122").Value = 1572
$ws.Range("L122").Value = 14148
$ws.Range("N122").Value = -19048

$ws.Range("H129").Value = 1860
$ws.Range("I129").Value = 1575
$ws.Range("K129").Value = 4725
$ws.Range("M129").Value = 275

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 17000
$ws.Range("J63").Value = 17000
$ws.Range("L63").Value = 17000
$ws.Range("N63").Value = -18248

$ws.Range("H66").Value = 17000
$ws.Range("J66").Value = 17000
$ws.Range("L66").Value = 51000
$ws.Range("N66").Value = -57240

$ws.Range("H132").Value = 2479.7
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3113.8572
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 9341.571599999999
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -14401.5716

$ws.Range("H136").Value = 1866.3334
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
